$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the B/C-quarter rows that got reordered (A column label + B:E data)
# for each of the four years present in the sheet.
$swapPairs = @(
    @(3, 4),
    @(7, 8),
    @(11, 12),
    @(15, 16)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("A$($r1):E$($r1)")
    $range2 = $ws.Range("A$($r2):E$($r2)")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Drop the F (production-sales ratio) and G (sales volume, non-cumulative)
# columns entirely - they were removed from the dataset.
$ws.Range("F1:G17").Delete()
